# Add a new "2021" data column (column R) to the worksheet, mirroring
# the formatting of the existing "2020" column (column Q), then fill in
# the new values for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles/number formats) from column Q (rows 4-34) into
# the new column R so the new column visually matches the rest of the table.
[void]$ws.Range("Q4:Q34").Copy($ws.Range("R4:R34"))

# New header value for the added year column.
$ws.Range("R4").Value = 2021

# New data values for the "2021" column.
$ws.Range("R5").Value = 0.8
$ws.Range("R6").Value = 0.4
$ws.Range("R7").Value = 1.2
$ws.Range("R8").Value = 0.2
$ws.Range("R9").Value = "-"
$ws.Range("R10").Value = 0.4
$ws.Range("R11").Value = 0.6
$ws.Range("R12").Value = 0.8
$ws.Range("R13").Value = 0.5
$ws.Range("R14").Value = 0.4
$ws.Range("R15").Value = "-"
$ws.Range("R16").Value = 0.8
$ws.Range("R17").Value = 0.3
$ws.Range("R18").Value = 0.7
$ws.Range("R19").Value = "-"
$ws.Range("R20").Value = 0.5
$ws.Range("R21").Value = 0.1
$ws.Range("R22").Value = 0.8
$ws.Range("R23").Value = 1.1
$ws.Range("R24").Value = 1.5
$ws.Range("R25").Value = 0.7
$ws.Range("R26").Value = 2.2
$ws.Range("R27").Value = 1
$ws.Range("R28").Value = 3.5
$ws.Range("R29").Value = 0.8
$ws.Range("R30").Value = 0.2
$ws.Range("R31").Value = 1.6
$ws.Range("R32").Value = 0.3
$ws.Range("R33").Value = "-"
$ws.Range("R34").Value = 0.6

# Update the active selection to match the author's last position.
[void]$ws.Range("Q11").Select()
